{"js": "// Replace the three-digit x one-digit multiplication prompts throughout\n// the document body with the new values from the commit.\nconst replacements = [\n  [\"179\u00d72=\", \"896\u00d72=\"],\n  [\"275\u00d77=\", \"215\u00d74=\"],\n  [\"483\u00d72=\", \"159\u00d74=\"],\n  [\"259\u00d74=\", \"430\u00d72=\"],\n  [\"867\u00d78=\", \"737\u00d76=\"],\n  [\"952\u00d72=\", \"216\u00d72=\"],\n  [\"766\u00d76=\", \"412\u00d79=\"],\n  [\"878\u00d77=\", \"835\u00d78=\"],\n  [\"617\u00d78=\", \"978\u00d76=\"],\n  [\"415\u00d77=\", \"516\u00d73=\"],\n  [\"494\u00d74=\", \"361\u00d73=\"],\n  [\"481\u00d79=\", \"811\u00d73=\"],\n  [\"998\u00d76=\", \"453\u00d72=\"],\n  [\"939\u00d77=\", \"161\u00d75=\"],\n  [\"608\u00d76=\", \"801\u00d72=\"],\n  [\"545\u00d73=\", \"985\u00d77=\"],\n  [\"935\u00d77=\", \"920\u00d79=\"],\n  [\"345\u00d72=\", \"363\u00d73=\"],\n  [\"382\u00d74=\", \"414\u00d77=\"],\n  [\"185\u00d79=\", \"224\u00d76=\"],\n  [\"620\u00d76=\", \"493\u00d74=\"],\n  [\"160\u00d75=\", \"914\u00d78=\"],\n  [\"257\u00d79=\", \"134\u00d79=\"],\n  [\"919\u00d72=\", \"263\u00d79=\"],\n  [\"188\u00d72=\", \"698\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"179\u00d72=\", \"896\u00d72=\"),\n    @(\"275\u00d77=\", \"215\u00d74=\"),\n    @(\"483\u00d72=\", \"159\u00d74=\"),\n    @(\"259\u00d74=\", \"430\u00d72=\"),\n    @(\"867\u00d78=\", \"737\u00d76=\"),\n    @(\"952\u00d72=\", \"216\u00d72=\"),\n    @(\"766\u00d76=\", \"412\u00d79=\"),\n    @(\"878\u00d77=\", \"835\u00d78=\"),\n    @(\"617\u00d78=\", \"978\u00d76=\"),\n    @(\"415\u00d77=\", \"516\u00d73=\"),\n    @(\"494\u00d74=\", \"361\u00d73=\"),\n    @(\"481\u00d79=\", \"811\u00d73=\"),\n    @(\"998\u00d76=\", \"453\u00d72=\"),\n    @(\"939\u00d77=\", \"161\u00d75=\"),\n    @(\"608\u00d76=\", \"801\u00d72=\"),\n    @(\"545\u00d73=\", \"985\u00d77=\"),\n    @(\"935\u00d77=\", \"920\u00d79=\"),\n    @(\"345\u00d72=\", \"363\u00d73=\"),\n    @(\"382\u00d74=\", \"414\u00d77=\"),\n    @(\"185\u00d79=\", \"224\u00d76=\"),\n    @(\"620\u00d76=\", \"493\u00d74=\"),\n    @(\"160\u00d75=\", \"914\u00d78=\"),\n    @(\"257\u00d79=\", \"134\u00d79=\"),\n    @(\"919\u00d72=\", \"263\u00d79=\"),\n    @(\"188\u00d72=\", \"698\u00d77=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $r = $d.Content\n    $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
